# Applies the "additional scraping" edit:
#  1. Inserts a new "Player Info" worksheet BEFORE the existing "ODI Batting"
#     sheet, with header row (ID, NAME, BATTING_HAND, BOWL_STYLE) and one
#     data row for player 7120 (Stephen Thomas Doheny).
#  2. On the "ODI Batting" sheet, renames column D's header from
#     MATCH_CARD_LINK to MATCH_CODE, and replaces each data row's full
#     howstat.com scorecard URL with just the bare match code number
#     (kept as text, not a numeric value).

$wb = $excel.ActiveWorkbook

# ---- existing sheet reference (by name, so it is stable) ----
$battingName = $wb.Worksheets.Item(1).Name

# ---- 1. add the new "Player Info" sheet before it ----
$beforeSheet = $wb.Worksheets.Item($battingName)
$info = $wb.Worksheets.Add($beforeSheet)
$info.Name = "Player Info"

# Header row
$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

$header = $info.Range("A1:D1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Match the page margins used elsewhere in this workbook (values in points).
$info.PageSetup.LeftMargin = 54
$info.PageSetup.RightMargin = 54
$info.PageSetup.TopMargin = 72
$info.PageSetup.BottomMargin = 72
$info.PageSetup.HeaderMargin = 36
$info.PageSetup.FooterMargin = 36

# Data row - ID must stay text (not get coerced to a number)
$info.Range("A2").NumberFormat = "@"
$info.Range("A2").Value = "7120"
$info.Range("A2").Style = "Normal"

$info.Range("B2").Value = "Stephen Thomas Doheny"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Does Not Bowl | Unknown"

$info.Range("A1").Select()

# ---- 2. update the "ODI Batting" sheet's MATCH_CARD_LINK column ----
# Re-fetch by name: inserting the sheet above shifted index-based refs.
$batting = $wb.Worksheets.Item($battingName)

$batting.Range("D1").Value = "MATCH_CODE"

$codes = @{
    2 = "4693"
    3 = "4694"
    4 = "4696"
    5 = "4726"
    6 = "4729"
    7 = "4734"
}

foreach ($row in $codes.Keys) {
    $cell = $batting.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $codes[$row]
    $cell.Style = "Normal"
}
